# Updated search command system test.
#
# 1) Bump the cached "datetimeFigureOut" field text (Insert > Header & Footer
#    "Date and time" placeholder) from 12/15/2018 to 4/7/2019 everywhere it is
#    stored: the slide master, every slide layout, and the notes master.
# 2) Remove the "Cloud 50" shape and its "Elbow Connector 51" connector from
#    slide 1.

$p = $ppt.ActivePresentation

$oldDate = "12/15/2018"
$newDate = "4/7/2019"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master's Date Placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's Date Placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master's Date Placeholder.
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes

# Remove the "Cloud 50" shape and its connector ("Elbow Connector 51") from
# slide 1. Delete the connector first so shape indices of the shapes we still
# need to look up don't shift underneath us.
$slide = $p.Slides.Item(1)
for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "Elbow Connector 51" -or $shp.Name -eq "Cloud 50") {
        $shp.Delete()
    }
}
